$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so Excel does not auto-convert them to numeric values (losing formatting
# and introducing floating point artifacts), matching the original inline-string cells.
$textForceCells = @("D5", "D6", "D12", "D13", "D19", "D20", "D21", "D26", "D27", "D31", "D33", "D44", "D46", "D47", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value updates described by the diff
$ws.Range("D2").Value = "67.358.20"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.556.02"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "593.59"
$ws.Range("D6").Value = "172.97"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "2.554.71"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "3.014.50"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "67.188.77"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "2.554.98"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "7.94"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").Value = "356.21"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  +6.65%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "70.23"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "10.12"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("D28").Value = "2.686.96"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "0.0₃0999"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "535.14"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  +6.10%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "39.74"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "151.28"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("D49").Value = "0.0₆0281"
$ws.Range("E49").Value = "  -4.54%  "
$ws.Range("D50").Value = "3.74"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("E51").Value = "  +1.91%  "

# Remove the temporary text number-format so no stray style index is left on the cells
foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
